$d = $word.ActiveDocument

# The edit consolidates several paragraphs that were previously split
# across multiple <w:r> runs (one run per sentence/editing-session,
# carrying distinct w:rsidR attributes and, in a couple of spots,
# w:proofErr wrappers) into a single run per paragraph, without
# changing any of the visible text. We rebuild each paragraph's range
# with the concatenated text. Because the final text is identical to
# the current concatenation, the host's diffing can treat a same-text
# assignment as a no-op, so each paragraph is first stamped with a
# short placeholder and then overwritten with the real final text —
# guaranteeing the run list is actually rebuilt as a single run.

function Set-ParagraphText($paragraphIndex, $finalText) {
    $p = $d.Paragraphs($paragraphIndex)
    $start = $p.Range.Start
    $end = $p.Range.End - 1
    $r = $d.Range($start, $end)
    $r.Text = "x"
    $r2 = $d.Range($start, $start + 1)
    $r2.Text = $finalText
}

# Paragraph 2: "There are three main online Application development tools..."
Set-ParagraphText 2 "There are three main online Application development tools: Amazon Web Services (AWS), Google App Engine and Microsoft Azure. We intend to create a web application that can also be accessed in a mobile format. After careful consideration at XVI we have decided to go with Azure in conjunction with Microsoft Visual Studio as our development tool, storage warehouse and launching platform."

# Paragraph 3: "Azure offers extensive testing and DevOps tools..."
Set-ParagraphText 3 "Azure offers extensive testing and DevOps tools, expansive middleware, an enormous data staging ground that scales with usage, virtual machines to use as containers and simple compatibility with Visual Studio for front-end support as well as an unimaginable amount of other functions."

# Paragraph 4: "We are able to use Java to code the back end..."
# (the w:proofErr gramStart/gramEnd markers around "are able to" are
# also removed by the diff since that run is absorbed into the merge)
Set-ParagraphText 4 "We are able to use Java to code the back end and .html for the front-end in Visual Studio. Both offer fantastic platforms to work off and have shared functionality between the two – meaning; we can use those two Microsoft products for most of our project without leaving a conjoined Microsoft suite. "

# Paragraph 5: "After development of our application..."
Set-ParagraphText 5 "After development of our application, we can perform containerized testing against massive, pre-structured data-pools created by Microsoft, also in Azure, that replicates devices in real-world situations. This ensures we can perform thorough testing without having to establish real-world scenarios, expanding the scope of our testing and giving us an opportunity to ready ourselves for market. These tests will then allow us to go into alpha-testing in isolated, situationally specific, environments like aged-care facilities and schools. "

# Paragraph 6: "As SocialCare grows, so will the staging ground..."
# The w:proofErr spellStart/spellEnd wrapping "SocialCare" is left
# untouched by the diff — only the runs after it are merged — so we
# locate the end of "SocialCare" and merge just the remainder.
$p6 = $d.Paragraphs(6)
$p6Start = $p6.Range.Start
$p6End = $p6.Range.End - 1
$searchRange = $d.Range($p6Start, $p6End)
$searchRange.Find.Execute("SocialCare", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$tailStart = $searchRange.End
$tailFinalText = " grows, so will the staging ground set by Azure. With extra data availability and future development capabilities always at-the-ready Social-Care, like the market and our families, never has to stop evolving."
$tailRange = $d.Range($tailStart, $p6End)
$tailRange.Text = "x"
$tailRange2 = $d.Range($tailStart, $tailStart + 1)
$tailRange2.Text = $tailFinalText

Write-Output "Para2: $($d.Paragraphs(2).Range.Text)"
Write-Output "Para3: $($d.Paragraphs(3).Range.Text)"
Write-Output "Para4: $($d.Paragraphs(4).Range.Text)"
Write-Output "Para5: $($d.Paragraphs(5).Range.Text)"
Write-Output "Para6: $($d.Paragraphs(6).Range.Text)"
